$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.985.67'
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").Value = '2.303.52'
$ws.Range("E3").Value = '  -0.34%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '''299.82'
$ws.Range("E5").Value = '  -0.72%  '

$ws.Range("D6").Value = '''97.45'
$ws.Range("E6").Value = '  -1.21%  '

$ws.Range("E7").Value = '  -1.71%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = '''0.506'
$ws.Range("E9").Value = '  -3.06%  '

$ws.Range("D10").Value = '''35.71'
$ws.Range("E10").Value = '  +0.09%  '

$ws.Range("D11").Value = '''0.0789'
$ws.Range("E11").Value = '  -0.29%  '

$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").Value = '''18.04'
$ws.Range("E12").Value = '  +0.63%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.118'
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("D14").Value = '''6.79'
$ws.Range("E14").Value = '  -1.53%  '

$ws.Range("D15").Value = '2.659.41'
$ws.Range("E15").Value = '  -0.46%  '

$ws.Range("D16").Value = '2.282.28'
$ws.Range("E16").Value = '  +0.80%  '

$ws.Range("D17").Value = '''0.778'
$ws.Range("E17").Value = '  -1.43%  '

$ws.Range("D18").Value = '42.905.06'
$ws.Range("E18").Value = '  -0.15%  '

$ws.Range("D19").Value = '''12.71'
$ws.Range("E19").Value = '  -6.03%  '

$ws.Range("D20").Value = '0.0₃0904'
$ws.Range("E20").Value = '  -0.67%  '

$ws.Range("E21").Value = '  -2.55%  '

$ws.Range("D22").Value = '''67.96'
$ws.Range("E22").Value = '  -0.56%  '

$ws.Range("D23").Value = '''240.26'
$ws.Range("E23").Value = '  +0.16%  '

$ws.Range("D24").Value = '''2.14'
$ws.Range("E24").Value = '  -1.40%  '

$ws.Range("E26").Value = '  -1.14%  '

$ws.Range("E27").Value = '  -0.32%  '

$ws.Range("D28").Value = '''25.49'
$ws.Range("E28").Value = '  +2.39%  '

$ws.Range("D29").Value = '''165.66'
$ws.Range("E29").Value = '  -1.48%  '

$ws.Range("D30").Value = '''2.02'
$ws.Range("E30").Value = '  -1.06%  '

$ws.Range("D31").Value = '''9.05'
$ws.Range("E31").Value = '  -1.43%  '

$ws.Range("D32").Value = '''33.12'
$ws.Range("E32").Value = '  -0.71%  '

$ws.Range("D33").Value = '''4.96'
$ws.Range("E33").Value = '  +2.06%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").Value = '''5.04'
$ws.Range("E35").Value = '  -3.14%  '

$ws.Range("D36").Value = '''16.97'
$ws.Range("E36").Value = '  -7.15%  '

$ws.Range("D37").Value = '''2.38'
$ws.Range("E37").Value = '  -1.30%  '

$ws.Range("D38").Value = '''0.0687'
$ws.Range("E38").Value = '  -0.99%  '

$ws.Range("E39").Value = '  -1.39%  '

$ws.Range("E40").Value = '  -2.23%  '

$ws.Range("E41").Value = '  -1.49%  '

$ws.Range("D42").Value = '''2.74'
$ws.Range("E42").Value = '  -0.97%  '

$ws.Range("D43").Value = '2.010.41'
$ws.Range("E43").Value = '  +0.50%  '

$ws.Range("E44").Value = '  -2.38%  '

$ws.Range("D45").Value = '''10.14'
$ws.Range("E45").Value = '  +0.27%  '

$ws.Range("E46").Value = '  -0.83%  '

$ws.Range("D47").Value = '''17.27'
$ws.Range("E47").Value = '  -1.88%  '

$ws.Range("E48").Value = '  -1.38%  '

$ws.Range("E49").Value = '  -4.24%  '

$ws.Range("D50").Value = '''53.60'
$ws.Range("E50").Value = '  -2.34%  '

$ws.Range("D51").Value = '2.527.89'
$ws.Range("E51").Value = '  -0.45%  '
